$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").EntireColumn.Delete()

$ws.Range("G15").Select()
